$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) and, after the shift, the "SC 92" row
# (which becomes row 27) so the remaining rows close the gap - matches
# the dimension shrinking from A1:F35 to A1:F33.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Update the F-column (imputed) values that changed between the two
# error-calculation passes.
$ws.Range("F2").Value = 18.03
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F11").Value = 17.65
$ws.Range("F13").ClearContents()
$ws.Range("F21").Value = 16.58
$ws.Range("F25").ClearContents()
$ws.Range("F33").Value = 17.53
